$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source data is scraped text (e.g. prices like "1.001"); force Text format
# on any price cell whose new value would otherwise be auto-parsed as a number,
# matching the original workbook convention of storing Price as text.

$ws.Range("D2").Value = "23.405.66"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.629.41"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.05"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3785"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3648"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.60"
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08226"
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.232"
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.34"
$ws.Range("E13").Value = "  -2.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.543"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001250"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.319"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "1.628.59"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.99"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06988"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.65"
$ws.Range("E20").Value = "  -2.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.499"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.71"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").Value = "23.413.17"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.122"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.451"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.39"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.06"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.296"
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.93"
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("D31").Value = "1.809.49"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.260"
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.822"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.020"
$ws.Range("E34").Value = "  +5.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.83"
$ws.Range("E35").Value = "  +5.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02782"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2520"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08780"
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.11"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.10"
$ws.Range("E44").Value = "  -3.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6541"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.302"
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.973"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.195"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.92"
$ws.Range("E51").Value = "  -2.09%  "

# Rows 39-42: Hedera/InternetComputer(DFINITY) and TrustWalletToken/TheSandbox swapped rank
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07095"
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.017"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.349"
$ws.Range("E41").Value = "  -2.18%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7022"
$ws.Range("E42").Value = "  -1.02%  "
